# Minor databook ordering fix
#
# "Workbook table parsing is currently rigidly in the order that table
#  objects are generated in the corresponding page specifications. ...
#  tables must be ordered appropriately."
#
# The "Parameters" sheet lists five small lookup tables, one after another
# (a header row + a data row, separated by a blank row): rows 1-2, 4-5,
# 7-8, 10-11 and 13-14. Each table's data row holds a "Quantity Type"
# (column B) and a "Constant" (column C) for the parameter named in
# column A of the header row.
#
# The tables were generated in the wrong order; this reorders the
# *content* of the five tables (label, quantity type and constant/formula)
# while leaving the table scaffolding (headers, years, styles) in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# New content for each table slot, taken from whichever table is moving
# into that slot:
#   slot 1 (rows  1-2)  <- old table at rows 10-11 (Transmission probability)
#   slot 2 (rows  4-5)  <- old table at rows 13-14 (Number of contacts)
#   slot 3 (rows  7-8)  <- old table at rows  4-5  (Average duration)
#   slot 4 (rows 10-11) <- old table at rows  7-8  (Death rate, infected)
#   slot 5 (rows 13-14) <- old table at rows  1-2  (Death rate, susceptible)

$labelRows = @(1, 4, 7, 10, 13)
$dataRows  = @(2, 5, 8, 11, 14)

$labels = @(
    "Transmission probability per contact",
    "Number of contacts annually",
    "Average duration of infections (years)",
    "Death rate for infected people",
    "Death rate for susceptible people"
)

$quantityTypes = @(
    "Probability",
    "Number",
    "Duration",
    "Probability",
    "Probability"
)

# Either a plain numeric constant, or an IF(SUMPRODUCT(...)) formula string
# (matching whichever form the source table used).
$constants = @(
    0.008,
    '=IF(SUMPRODUCT(--(E5:W5<>""))=0,80,"N.A.")',
    5,
    '=IF(SUMPRODUCT(--(E11:W11<>""))=0,0.016,"N.A.")',
    '=IF(SUMPRODUCT(--(E14:W14<>""))=0,0.008,"N.A.")'
)

for ($i = 0; $i -lt 5; $i++) {
    $labelRow = $labelRows[$i]
    $dataRow  = $dataRows[$i]

    $ws.Cells.Item($labelRow, 1).Value = $labels[$i]
    $ws.Cells.Item($dataRow, 2).Value = $quantityTypes[$i]

    $cCell = $ws.Cells.Item($dataRow, 3)
    $val = $constants[$i]
    if ($val -is [string]) {
        $cCell.Formula = $val
    } else {
        $cCell.Value = $val
    }
}

# Data validation lists follow each row's new quantity type:
#   B2, B11, B14 -> Probability ; B8 -> Duration ; B5 -> Number
$dvCells = @("B2", "B5", "B8", "B11", "B14")
foreach ($addr in $dvCells) {
    $ws.Range($addr).Validation.Delete()
}

$dvSpecs = @{
    "B2"  = "Probability"
    "B11" = "Probability"
    "B14" = "Probability"
    "B8"  = "Duration"
    "B5"  = "Number"
}
foreach ($addr in @("B14", "B2", "B11", "B8", "B5")) {
    $listName = $dvSpecs[$addr]
    $ws.Range($addr).Validation.Add(3, 1, 1, ('"' + $listName + '"'))
}

# View state: "State Variables" keeps its remembered selection but is no
# longer the active tab; "Parameters" becomes active with rows 1:3 selected.
$stateVars = $wb.Worksheets.Item("State Variables")
$stateVars.Range("AC6").Select()

$params = $wb.Worksheets.Item("Parameters")
$params.Activate()
$params.Range("A1:XFD3").Select()
